# Update fixtures for Liga Suecia 2025: add new possession(%) / metadata
# columns, backfill first/second-half goal splits for several already
# recorded matches, and append the newly played fixtures (rows 146-153).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New header columns V:Y
# ---------------------------------------------------------------------------
$ws.Range("V1").Value = "Posesión Local (%)"
$ws.Range("W1").Value = "Posesión Visita (%)"
$ws.Range("X1").Value = "fuente_tiempos"
$ws.Range("Y1").Value = "estado_datos"

# ---------------------------------------------------------------------------
# 2) Corrections to existing rows (half-time goal splits / possession)
# ---------------------------------------------------------------------------
$corrections = @(
    ,("Q117", 49)
    ,("R117", 51)
    ,("M119", 1)
    ,("O119", 2)
    ,("M120", 1)
    ,("O120", 2)
    ,("M122", 3)
    ,("O122", 1)
    ,("N123", 2)
    ,("P123", 0)
    ,("M124", 1)
    ,("O124", 0)
    ,("N125", 1)
    ,("P125", 1)
    ,("N126", 1)
    ,("P126", 2)
    ,("N129", 2)
    ,("P129", 0)
    ,("Q130", 32)
    ,("R130", 68)
    ,("M131", 1)
    ,("O131", 2)
    ,("M132", 1)
    ,("N132", 1)
    ,("O132", 1)
    ,("P132", 2)
    ,("M133", 2)
    ,("O133", 1)
    ,("N135", 4)
    ,("P135", 2)
    ,("I137", 3)
    ,("J137", 2)
    ,("M137", 2)
    ,("N137", 1)
    ,("O137", 0)
    ,("P137", 2)
    ,("N138", 1)
    ,("P138", 0)
    ,("M139", 1)
    ,("N139", 3)
    ,("O139", 3)
    ,("P139", 0)
    ,("M140", 3)
    ,("N140", 4)
    ,("O140", 3)
    ,("P140", 0)
    ,("M141", 1)
    ,("N141", 1)
    ,("O141", 1)
    ,("P141", 1)
    ,("M142", 1)
    ,("N142", 1)
    ,("O142", 0)
    ,("P142", 0)
    ,("M145", 2)
    ,("O145", 1)
)

foreach ($pair in $corrections) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# Row 138's yellow-card counts (I138/J138) became unknown data -> blank them.
$ws.Range("I138").ClearContents()
$ws.Range("J138").ClearContents()

# ---------------------------------------------------------------------------
# 3) Append the new fixtures played 2025-08-09 .. 2025-08-11 (rows 146-153)
# ---------------------------------------------------------------------------
$newRows = @(
    ,@("2025-08-09","Malmo FF","Mjallby AIF",1,3,1342094,9,3,2,3,0,0,0,0,1,3,54,46,"V")
    ,@("2025-08-09","Halmstad","Sirius",0,1,1342096,8,8,3,1,0,0,0,0,0,1,50,50,"V")
    ,@("2025-08-10","AIK Stockholm","Djurgardens IF",0,0,1342090,3,5,2,2,0,0,0,0,0,0,43,57,"E")
    ,@("2025-08-10","Degerfors IF","BK Hacken",0,0,1342097,4,7,3,2,0,0,0,0,0,0,31,69,"E")
    ,@("2025-08-10","IFK Norrkoping","Hammarby FF",0,2,1342095,9,10,5,1,0,0,0,0,0,2,37,63,"V")
    ,@("2025-08-10","Osters IF","IF Brommapojkarna",1,1,1342093,3,1,1,5,0,0,1,0,0,1,39,61,"E")
    ,@("2025-08-11","IF Elfsborg","IFK Varnamo",2,2,1342091,8,3,1,0,0,0,1,0,1,2,54,46,"E")
    ,@("2025-08-11","Gais","IFK Goteborg",0,1,1342092,12,2,4,4,0,0,0,1,0,0,62,38,"V")
)

$r = 146
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value  = $row[0]   # Fecha
    $ws.Cells.Item($r, 2).Value  = $row[1]   # Local
    $ws.Cells.Item($r, 3).Value  = $row[2]   # Visita
    $ws.Cells.Item($r, 4).Value  = $row[3]   # Goles Local
    $ws.Cells.Item($r, 5).Value  = $row[4]   # Goles Visita
    $ws.Cells.Item($r, 6).Value  = $row[5]   # Fixture ID
    $ws.Cells.Item($r, 7).Value  = $row[6]   # Corners Local
    $ws.Cells.Item($r, 8).Value  = $row[7]   # Corners Visita
    $ws.Cells.Item($r, 9).Value  = $row[8]   # Amarillas Local
    $ws.Cells.Item($r, 10).Value = $row[9]   # Amarillas Visita
    $ws.Cells.Item($r, 11).Value = $row[10]  # Rojas Local
    $ws.Cells.Item($r, 12).Value = $row[11]  # Rojas Visita
    $ws.Cells.Item($r, 13).Value = $row[12]  # Goles 1T Local
    $ws.Cells.Item($r, 14).Value = $row[13]  # Goles 1T Visita
    $ws.Cells.Item($r, 15).Value = $row[14]  # Goles 2T Local
    $ws.Cells.Item($r, 16).Value = $row[15]  # Goles 2T Visita
    $ws.Cells.Item($r, 17).Value = $row[16]  # Posesión Local ()
    $ws.Cells.Item($r, 18).Value = $row[17]  # Posesión Visita ()
    $ws.Cells.Item($r, 19).Value = $row[18]  # Resultado
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 4) Style the new header cells like the rest of row 1 (bold, centered,
#    bordered) by copying the format already applied to column U's header.
# ---------------------------------------------------------------------------
$ws.Range("U1").Copy()
$ws.Range("V1:Y1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("V1").Value = "Posesión Local (%)"
$ws.Range("W1").Value = "Posesión Visita (%)"
$ws.Range("X1").Value = "fuente_tiempos"
$ws.Range("Y1").Value = "estado_datos"
